$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$t = $sh.Table
$t.Style.Id = "{C33BE8D6-3AA6-459B-BBBE-42C1B2030173}"
Write-Output "Style: $($t.Style)"
Write-Output "StyleId: $($t.Style.Id)"
